$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template/example row (row 5) held placeholder demo data ("1" in A5 and
# the "OO" placeholder text across B5:K5). Clear it out so the row is blank,
# ready for real data (matching the "(表七)台水公司災害通報彙整" sheet edit).
$ws.Range("A5:K5").ClearContents()

# Move the active selection from the old stray cell (Q13) to K7.
$ws.Range("K7").Select()
